# Apply the edits described by the diff to the workbook.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "sounds-meta-data"
$ws2 = $wb.Worksheets.Item(2)   # "sounds_list"

# ---------------------------------------------------------------
# Sheet1 ("sounds-meta-data") cell content updates.
# New shared strings are appended to the shared string table in the
# order they are first assigned, so the order below is chosen to
# reproduce the target shared string table ordering.
# ---------------------------------------------------------------

# Row 6 - was fire-alarm -> alarm, add All alert sounds (D6)
$ws1.Range("C6").Value = "alarm"
$ws1.Range("D6").Value = "All alert sounds"

# Row 8 - was Line-App-Official-Ring-Sounds -> ringtone, add D8
$ws1.Range("C8").Value = "ringtone"
$ws1.Range("D8").Value = "Line and WeChat default"

# Row 11 - was silence/98 -> pain-sound/57
$ws1.Range("B11").Value = 57
$ws1.Range("C11").Value = "pain-sound"

# Row 12 - was other-sounds/99 -> foot-setp/58
$ws1.Range("B12").Value = 58
$ws1.Range("C12").Value = "foot-setp"

# Row 10 - sneezing: add D10
$ws1.Range("D10").Value = "not-need now"

# Row 1 header - add new column E header
$ws1.Range("E1").Value = "first-sage included"

# Row 3/4/5/6 - "*" legend marker column E
$ws1.Range("E3").Value = "*"
$ws1.Range("E4").Value = "*"
$ws1.Range("E5").Value = "*"
$ws1.Range("E6").Value = "*"

# Row 3 - legend explanation text, larger font (16pt), centered, bordered
$ws1.Range("F3").Value = "*:表示用更多聲音進行訓練"
$ws1.Range("F3").Font.Size = 16
$ws1.Range("F3").HorizontalAlignment = -4108
$ws1.Range("F3").VerticalAlignment = -4108
$ws1.Range("F3").Borders.LineStyle = 1

# Row 5 - kettle-sound: add D5 (Chinese label)
$ws1.Range("D5").Value = "熱水壺笛聲"

# Row 13 - new row: silence/98 (reuses existing shared string)
$ws1.Range("A13").Value = 12
$ws1.Range("B13").Value = 98
$ws1.Range("C13").Value = "silence"

# Row 14 - new row: other-sounds/99 (reuses existing shared string)
$ws1.Range("A14").Value = 13
$ws1.Range("B14").Value = 99
$ws1.Range("C14").Value = "other-sounds"

# ---------------------------------------------------------------
# Sheet1 cosmetic/layout updates
# ---------------------------------------------------------------

# Column widths (approximate; engine quantizes to character-width units)
$ws1.Columns.Item(1).ColumnWidth = 16.71
$ws1.Columns.Item(2).ColumnWidth = 16
$ws1.Columns.Item(3).ColumnWidth = 43
$ws1.Columns.Item(4).ColumnWidth = 34.14
$ws1.Columns.Item(6).ColumnWidth = 41.14

# Print/page setup
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# Selection moves to D6
[void]$ws1.Range("D6").Select()

# ---------------------------------------------------------------
# Sheet2 ("sounds_list") - no content changes needed; the shared
# string table / indices are recomputed automatically on save.
# ---------------------------------------------------------------

Write-Output "edits applied"
